# chore: add pfp file names
#
# Adds profile-picture file names to the "team" sheet's pfp_file_name column (D),
# and updates Leonardo Bosche's research_info bio text.

$wb = $excel.ActiveWorkbook
$wsTeam = $wb.Worksheets.Item("team")
$wsPubs = $wb.Worksheets.Item("publications")

# Map each row (by name in column B) to its profile picture file name (column D)
$pfpByName = @{
    "Ignacio Ciampitti" = "ignacio_ciampitti.jpg"
    "Pedro Cisdeli"     = "pedro_cisdeli.png"
    "Gustavo Santiago"  = "gustavo_santiago.jpg"
    "Leonardo Bosche"   = "leonardo_bosche.jpg"
    "Natalia Volpato"   = "natalia_volpato.jpg"
    "Federico Gomez"    = "federico_gomez.jpg"
    "German Mandrini"   = "german_mandrini.png"
    "Priscila Cano"     = "priscila_cano.jpg"
    "Franco Muriñigo"   = "franco_murinigo.jpg"
    "Thatiane Pereira"  = "thatiane_pereira.jpg"
}

for ($row = 2; $row -le 11; $row++) {
    $name = $wsTeam.Cells.Item($row, 2).Value()
    $pfp = $pfpByName[$name]
    if ($pfp) {
        $wsTeam.Cells.Item($row, 4).Value = $pfp
    }
}

# Update Leonardo Bosche's research_info bio (row 5, column F)
$wsTeam.Range("F5").Value = "Agronomist developing decision support models for precision nitrogen management in corn. Integrates crop physiology, sensing technologies, and statistics to improve productivity while reducing environmental impact."

# Restore cursor/selection state left by the author while editing
[void]$wsPubs.Range("D19").Select()
[void]$wsTeam.Range("F5").Select()
$wsTeam.Activate()
